# Applies the cryptos.xlsx price/volume update (Fri Aug 4 21:45:53 UTC 2023 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.015.86"
$ws.Range("E2").Value = "  -0.93%  "

# Row 3
$ws.Range("D3").Value = "1.826.44"
$ws.Range("E3").Value = "  -0.90%  "

# Row 4
$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'241.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
$ws.Range("D6").Value = "'0.6311"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.06%  "

# Row 8
$ws.Range("D8").Value = "'44.59"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.82%  "

# Row 9
$ws.Range("D9").Value = "'0.2929"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.32%  "

# Row 10
$ws.Range("D10").Value = "'0.07321"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.56%  "

# Row 11
$ws.Range("D11").Value = "'22.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.23%  "

# Row 12
$ws.Range("D12").Value = "'0.07669"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.65%  "

# Row 13
$ws.Range("D13").Value = "1.829.46"
$ws.Range("E13").Value = "  -0.39%  "

# Row 14
$ws.Range("D14").Value = "'4.983"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.52%  "

# Row 15
$ws.Range("D15").Value = "'0.6644"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "

# Row 16
$ws.Range("D16").Value = "'82.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.67%  "

# Row 17
$ws.Range("D17").Value = "'6.066"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.43%  "

# Row 18
$ws.Range("D18").Value = "'0.000008645"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.03%  "

# Row 19
$ws.Range("D19").Value = "28.945.68"
$ws.Range("E19").Value = "  -1.08%  "

# Row 20
$ws.Range("D20").Value = "2.081.09"
$ws.Range("E20").Value = "  +0.14%  "

# Row 21
$ws.Range("D21").Value = "'12.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.04%  "

# Row 22
$ws.Range("D22").Value = "'223.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.93%  "

# Row 23
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").Value = "'7.157"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "

# Row 26
$ws.Range("D26").Value = "'158.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.84%  "

# Row 27
$ws.Range("D27").Value = "'8.441"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.09%  "

# Row 28
$ws.Range("D28").Value = "'0.1372"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.32%  "

# Row 29
$ws.Range("D29").Value = "'17.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.81%  "

# Row 30
$ws.Range("D30").Value = "'1.507"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.54%  "

# Row 31
$ws.Range("D31").Value = "'4.092"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.63%  "

# Row 32
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.201"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.62%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.016"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.30%  "

# Row 34
$ws.Range("D34").Value = "'0.05296"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "

# Row 35
$ws.Range("D35").Value = "'0.7408"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.61%  "

# Row 36
$ws.Range("D36").Value = "'1.826"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.91%  "

# Row 37
$ws.Range("D37").Value = "'1.150"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.09%  "

# Row 38
$ws.Range("D38").Value = "'2.653"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.04%  "

# Row 39
$ws.Range("D39").Value = "1.291.10"

# Row 40
$ws.Range("D40").Value = "'2.748"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.69%  "

# Row 41
$ws.Range("D41").Value = "'0.01780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.51%  "

# Row 42
$ws.Range("D42").Value = "'6.334"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.03%  "

# Row 43
$ws.Range("D43").Value = "'0.8904"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.07%  "

# Row 44
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.82%  "

# Row 45
$ws.Range("D45").Value = "'102.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "

# Row 46
$ws.Range("D46").Value = "1.979.01"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("D47").Value = "'0.5139"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.46%  "

# Row 48
$ws.Range("D48").Value = "'64.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "

# Row 49
$ws.Range("D49").Value = "'0.00000000119"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.87%  "

# Row 50
$ws.Range("D50").Value = "'1.729"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.70%  "

# Row 51
$ws.Range("D51").Value = "'0.05828"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.97%  "
